# Update the standard-error / variance table in the testing sheet with the
# new, more sensible simulation parameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "j" coefficient table (rows 3-7) -------------------------------------
$ws.Range("B3").Value = "(0.222)"
$ws.Range("C3").Value = "(0.250)"

$ws.Range("B4").Value = "(0.250)"
$ws.Range("C4").Value = "(0.287)"

$ws.Range("B5").Value = "(0.294)"
$ws.Range("C5").Value = "(0.331)"

$ws.Range("B6").Value = "(0.347)"
$ws.Range("C6").Value = "(0.382)"

$ws.Range("B7").Value = "(0.402)"
$ws.Range("C7").Value = "(0.433)"

# --- "x" row (row 9) -------------------------------------------------------
$ws.Range("B9").Value = "(0.288)"
$ws.Range("C9").Value = "(0.285)"

# --- variance components (rows 13-15) --------------------------------------
$ws.Range("B13").Value = "(0.491)"
$ws.Range("C13").Value = "(0.546)"

$ws.Range("B14").Value = "(1.857)"
$ws.Range("C14").Value = "(2.371)"

$ws.Range("B15").Value = "(0.960)"
$ws.Range("C15").Value = "(0.946)"

# --- second-model-only standard errors (rows 17-21, column C) --------------
$ws.Range("C17").Value = "(0.061)"
$ws.Range("C18").Value = "(0.051)"
$ws.Range("C19").Value = "(0.006)"
$ws.Range("C20").Value = "(0.039)"
$ws.Range("C21").Value = "(0.025)"

# --- number of observations (row 22) ----------------------------------------
$ws.Range("B22").Value = 8050
$ws.Range("C22").Value = 9932
